# The "Primer Apellido" (first surname) column is no longer a field the
# application asks for, so its data (header + values) is removed from the
# sheet while keeping every other column (B, D, E, F, G, H) exactly where
# it is - i.e. column C becomes empty, nothing shifts.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$colC = $ws.Range("C1:C6")
$colC.ClearContents() | Out-Null

# Leave the selection on the column that was just cleared, matching the
# state the workbook was saved in.
$colC.Select() | Out-Null
